# Weekly refresh: insert the newest "Ajo" price record at the top of the
# data block (row 21) for "Vega Monumental Concepción", pushing the
# existing history (previously rows 21-124) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 21; everything below (old rows 21-124) shifts
# down to 22-125, carrying its data and formatting with it.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with this week's record.
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44558
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 100112003
$ws.Range("G21").Value = "Ajo"
$ws.Range("H21").Value = "Chino"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 16000
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 16400
$ws.Range("N21").Value = "$/caja 10 kilos"
$ws.Range("O21").Value = "China"
$ws.Range("P21").Value = 1640
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = "Hortaliza"
